$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update floodmedia column (H) for data rows 2-27 from numeric 0 to text "None"
$ws.Range("H2:H27").Value = "None"

# Update the active selection to reflect the edited range
$ws.Range("H2:H27").Select()
